# ---------------------------------------------------------------------------
# Applies the "Updates to transmission capacity and GDP files" edit to the
# BAU GDP workbook:
#   - OECD Data: rename G1 header, add H (National Growth) and I (State GDP)
#     columns with growth-rate / compounded-state-GDP formulas.
#   - BGDP: re-point the GDP formulas at the new OECD Data!I column instead
#     of recomputing from the national conversion factors.
#   - Re-create the view/selection state (active sheet, selections) that
#     Excel records when a user was last working on the BGDP sheet with the
#     OECD Data sheet scrolled/selected at K6.
# ---------------------------------------------------------------------------

$wb  = $excel.ActiveWorkbook
$wsAbout = $wb.Worksheets.Item("About")
$wsOecd  = $wb.Worksheets.Item("OECD Data")
$wsBgdp  = $wb.Worksheets.Item("BGDP")

# ---------------------------------------------------------------------------
# 1. OECD Data sheet: headers
# ---------------------------------------------------------------------------
$wsOecd.Range("G1").Value = "National GDP"
$wsOecd.Range("H1").Value = "National Growth"
$wsOecd.Range("I1").Value = "State GDP"

# ---------------------------------------------------------------------------
# 2. OECD Data sheet: column H growth-rate formulas (row 2 .. 48)
#    H2 = G2/G2 (first year: growth of 1), Hn = Gn/G(n-1) afterwards.
# ---------------------------------------------------------------------------
$wsOecd.Range("H2").Formula = "=G2/G2"
for ($r = 3; $r -le 48; $r++) {
    $prev = $r - 1
    $wsOecd.Cells.Item($r, 8).Formula = "=G$r/G$prev"
}
# H3 would otherwise inherit G3's custom (2-decimal) number style when the
# cell is created next to it; reset it back to the default/Normal style so
# it matches the rest of the freshly-created column H cells.
$wsOecd.Range("H3").Style = "Normal"

# ---------------------------------------------------------------------------
# 3. OECD Data sheet: column I compounded State GDP formulas (row 8 .. 48)
#    Rows 2-7 keep their existing hard-coded historical values.
#    In = I(n-1) * Hn starting at row 8.
# ---------------------------------------------------------------------------
for ($r = 8; $r -le 48; $r++) {
    $prev = $r - 1
    $wsOecd.Cells.Item($r, 9).Formula = "=I$prev*H$r"
}

# ---------------------------------------------------------------------------
# 4. OECD Data sheet: column width touch-ups (stored OOXML width is
#    ColumnWidth + 5/6 in this engine, so we back-compute the ColumnWidth
#    to hit the exact target stored widths).
# ---------------------------------------------------------------------------
$wsOecd.Columns.Item(1).ColumnWidth = 12.330729166666666   # -> 13.1640625
$wsOecd.Columns.Item(3).ColumnWidth = 9.998697916666666    # -> 10.83203125
$wsOecd.Columns.Item(5).ColumnWidth = 12.498697916666666   # -> 13.33203125
$wsOecd.Columns.Item(7).ColumnWidth = 16.330729166666668   # -> 17.1640625
$wsOecd.Columns.Item(8).ColumnWidth = 12.998697916666666   # -> 13.83203125
$wsOecd.Columns.Item(9).ColumnWidth = 10.998697916666666   # -> 11.83203125

# ---------------------------------------------------------------------------
# 5. BGDP sheet: re-point GDP formulas at OECD Data!I instead of recomputing
#    from About!$A$15 / About!$A$16 conversion factors.
# ---------------------------------------------------------------------------
for ($r = 2; $r -le 48; $r++) {
    $wsBgdp.Cells.Item($r, 2).Formula = "='OECD Data'!I$r"
}

# ---------------------------------------------------------------------------
# 6. BGDP sheet: column width touch-ups.
# ---------------------------------------------------------------------------
$wsBgdp.Columns.Item(1).ColumnWidth = 16.830729166666668   # -> 17.6640625
$wsBgdp.Columns.Item(2).ColumnWidth = 14.330729166666666   # -> 15.1640625

# ---------------------------------------------------------------------------
# 7. About sheet: column width touch-up.
# ---------------------------------------------------------------------------
$wsAbout.Columns.Item(2).ColumnWidth = 76.66666666666667   # -> 77.5

# ---------------------------------------------------------------------------
# 8. Recreate the recorded view/selection state: OECD Data selection at K6,
#    then BGDP active/selected at E43 (BGDP ends up the active sheet/tab).
# ---------------------------------------------------------------------------
$wsOecd.Activate()
$wsOecd.Range("K6").Select()

$wsBgdp.Activate()
$wsBgdp.Range("E43").Select()

Write-Output "BAU GDP workbook updated"
